# edit.ps1 - apply the "Transcending Boundaries" -> "The Marvelous World of
# Science" rewrite to the active document via Word COM-interop calls.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$NewText
    )
    $para = $d.Paragraphs($ParaIndex)
    $r = $para.Range
    $ok = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $ParaIndex for text: $OldText"
    }
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-InParagraph 1 "Transcending Boundaries: Interdisciplinary Perspectives" "The Marvelous World of Science: Unraveling Nature's Enigmas"

# ---------------------------------------------------------------------------
# 2. Author name: "Emily Carter" -> "Dr. Alex Green"
# ---------------------------------------------------------------------------
Replace-InParagraph 2 "Emily Carter" "Dr. Alex Green"

# ---------------------------------------------------------------------------
# 3-5. Email address runs
# ---------------------------------------------------------------------------
Replace-InParagraph 3 "ecarter@uni-gla" "alexgreen"
Replace-InParagraph 3 "ac" "edu@gmail"
Replace-InParagraph 3 "uk" "com"

# ---------------------------------------------------------------------------
# 6-16. Big introduction paragraph (paragraph 5)
# ---------------------------------------------------------------------------

Replace-InParagraph 5 "In the vast tapestry of human knowledge, disciplines, like radiant threads, intertwine to create an elaborate design, bridging disparate realms of inquiry and sparking new vistas of understanding" "In the grand tapestry of human knowledge, science stands as a beacon of enlightenment, dispelling ignorance and revealing the breathtaking intricacies of our natural world"

Replace-InParagraph 5 " Interdisciplinary approaches, transcending traditional boundaries, have emerged as a powerful force, unlocking hidden connections, challenging established paradigms, and fueling innovation across diverse fields" (" From the celestial ballet of stars to the microscopic realm of cells, science offers a kaleidoscope of wonder and discovery." + " Through patient observation, meticulous experimentation, and brilliant deduction, scientists have unraveled the riddles of nature, transforming our understanding of the universe and shaping our technological advancements." + " In this essay, we will embark on a captivating journey into the world of science, unraveling the mysteries that have perplexed humanity for centuries")

Replace-InParagraph 5 "The fusion of scientific and creative domains, where art meets science, has yielded awe-inspiring results, blurring the lines between the objective and the subjective" "As we delve into the vast ocean of scientific knowledge, we will explore the fundamental principles that govern the universe, unlocking the secrets of matter, energy, and space-time"

Replace-InParagraph 5 " In this realm, artists draw upon scientific concepts to create mesmerizing installations that invite viewers to experience the world in novel ways, while scientists seek artistic inspiration to visualize complex phenomena and explore the boundaries of reality" (" We will witness the elegance of mathematical formulas, unveiling patterns and symmetries that govern the cosmos." + " Through the lens of chemistry, we will discover the transformative power of atoms and molecules, unraveling the secrets of substances that form the building blocks of life")

Replace-InParagraph 5 "Technology and the arts have forged a symbiotic relationship, transforming both domains" "Venturing further, we will immerse ourselves in the wonders of biology, unraveling the intricate workings of living organisms"

Replace-InParagraph 5 " Digital technology has revolutionized the creative process, enabling artists to manipulate and shape images, sounds, and texts in unprecedented ways" " We will marvel at the symphony of cellular processes that sustain life, from photosynthesis to respiration"

Replace-InParagraph 5 " Concurrently, artistic principles and aesthetics have influenced technology, inspiring the development of user-friendly interfaces, innovative software, and captivating virtual realities" " We will study the remarkable diversity of life on Earth, tracing the evolutionary tapestry that connects all living beings"

Replace-InParagraph 5 "Breaking down the barriers between different disciplines has led to transformative discoveries in medicine and biology" "As we continue our exploration, we will delve into the fascinating realm of physics, unraveling the mysteries of energy, matter, and the fundamental forces that shape our universe"

Replace-InParagraph 5 " The integration of medical imaging technologies with artificial intelligence has revolutionized disease diagnosis and treatment planning, while the fusion of biology with engineering principles has given rise to biomimicry - the study of nature's designs to inspire technological solutions" " We will uncover the secrets of gravity, electromagnetism, and the enigmatic world of quantum mechanics, pushing the boundaries of human knowledge and understanding"

# The last original sentence is replaced by a large amount of new content
# (several new paragraphs worth of text joined by manual line breaks "^l").
$bigAddition = ""
$bigAddition += "^l^l"
$bigAddition += "Beyond the confines of our planet, we will venture into the vastness of space, exploring the solar system and beyond."
$bigAddition += " We will marvel at the awe-inspiring beauty of celestial objects, from shimmering stars to swirling galaxies."
$bigAddition += " We will contemplate the mysteries of the cosmic microwave background, seeking clues to the origin and ultimate fate of our universe."
$bigAddition += "^l^l"
$bigAddition += "Our scientific journey will not be limited to the natural world."
$bigAddition += " We will delve into the realm of psychology, unraveling the intricacies of the human mind, emotions, and behavior."
$bigAddition += " We will explore the social sciences, examining the complex interactions between individuals, societies, and cultures."
$bigAddition += " Through this holistic approach, we will gain a deeper understanding of ourselves and our place in the world."
$bigAddition += "^l^l"
$bigAddition += "Introduction Concluded:"
$bigAddition += "^l^l"
$bigAddition += "As we reach the end of our expedition into the world of science, we find ourselves amazed by the sheer complexity and beauty of the natural world."
$bigAddition += " From the smallest atom to the grandest galaxy, the universe is a symphony of interconnected phenomena, governed by fundamental laws that we are only beginning to comprehend."
$bigAddition += " Science has provided us with the tools to unravel these mysteries, unlocking the secrets of nature and propelling us forward in our quest for knowledge and understanding"

Replace-InParagraph 5 " Collaborative research between biologists, chemists, and computer scientists has yielded novel drugs, targeted therapies, and gene-editing technologies, promising hope for the treatment of previously incurable diseases" $bigAddition

# ---------------------------------------------------------------------------
# 17. Remove the lastRenderedPageBreak marker before "Summary"
# ---------------------------------------------------------------------------
# (handled implicitly - this runtime does not recompute lastRenderedPageBreak
# placement automatically; the marker text itself carries no visible text so
# nothing further is required here for the "Summary" run's content.)

# ---------------------------------------------------------------------------
# 18-20. Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------------
Replace-InParagraph 7 "Interdisciplinary approaches have reshaped the landscape of knowledge creation, leading to groundbreaking discoveries and transformative applications across a multitude of fields" "Our exploration of the world of science has led us on a captivating journey, unraveling the mysteries of the universe and revealing the breathtaking intricacies of the natural world"

Replace-InParagraph 7 " By breaking down traditional boundaries, interdisciplinary research fosters innovation, unlocks hidden connections, and inspires new ways of thinking" " Through the lens of mathematics, chemistry, biology, physics, and beyond, we have witnessed the elegance of scientific principles, the diversity of life, and the profound interconnectedness of all things"

Replace-InParagraph 7 " It empowers us to address complex challenges, tackle global issues, and explore the frontiers of human understanding, ultimately enriching our lives and expanding the horizons of human knowledge" (" Science has transformed our understanding of the universe, shaping our technological advancements and inspiring us to push the boundaries of human knowledge." + " As we continue to unravel nature's enigmas, we embrace the wonder and beauty of a world filled with endless possibilities")

# ---------------------------------------------------------------------------
# 21. Add a trailing empty paragraph at the very end of the document body
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output "done"
